$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

function Set-CellText($rowIndex, $newText) {
    $cell = $tbl.Cell($rowIndex, 1)
    $cell.Range.Text = $newText
}

# Single-value rows that simply change in place
Set-CellText 1  "0M"
Set-CellText 2  "0M"
Set-CellText 3  "0M"
Set-CellText 4  "400"
Set-CellText 6  "0.00053"
Set-CellText 7  "0.00022"
Set-CellText 8  "0.00004"
Set-CellText 9  "0.00037"
Set-CellText 10 "0.00042"
Set-CellText 11 "0.00044"
Set-CellText 12 "0.08734"

# Rows that collapse a tab-separated multi-run line into a single value
Set-CellText 44 "99.57"
Set-CellText 45 "0.09"
Set-CellText 46 "20"
